# Generate Report for Archive
#
# This applies two related changes to the localization-status workbook:
#   1. The status text "Ready for handoff" becomes "In Translation" on every
#      sheet that reports it (Overview!E2/F2, zh-cn!C2, de-de!C2).
#   2. Because the new status text is shorter than the old one, the status
#      columns that were sized to fit it are narrowed to match (Overview
#      columns E & F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# 1) Update the status value wherever it appears.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# 2) Narrow the columns that held the (now shorter) status text.
#    (12.5 is the ColumnWidth input that this engine's char-width rounding
#    lands closest to the target rendered width for these columns.)
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
